$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 10046.790000000001
$ws.Range("B9").Value = 9783.61
$ws.Range("C9").Value = 19.36
$ws.Range("D9").Value = 18.84
$ws.Range("E9").Value = $true
$ws.Range("F9").Value = -2.69
$ws.Range("G9").Value = 42612.673009259262
$ws.Range("H9").Value = $true
